# Apply "Add data for 2022-07-29" commit:
#  - Roll the "through July 20" snapshot forward to "through July 21"
#    (sheet name + running-month column header)
#  - Update the carjacking counts for the running month (column B) and
#    a handful of historical month columns across various neighborhoods

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet & update the column-B header label --------------------
$ws.Name = "Through 2022-07-21"
$ws.Range("B1").Value = "July 2022 (through July 21)"

# --- Cell value updates ---------------------------------------------------
# Row 3  - Englewood
$ws.Range("P3").Value = 4

# Row 4  - Auburn Gresham
$ws.Range("B4").Value = 4

# Row 5  - Garfield Park
$ws.Range("B5").Value = 8
$ws.Range("P5").Value = 11

# Row 8  - North Lawndale
$ws.Range("P8").Value = 15

# Row 20 - South Shore
$ws.Range("AK20").Value = 2

# Row 22 - Bridgeport
$ws.Range("I22").Value = 2

# Row 24 - South Deering
$ws.Range("W24").Value = 1
$ws.Range("AR24").Value = 2

# Row 29 - Humboldt Park
$ws.Range("B29").Value = 6
$ws.Range("I29").Value = 4

# Row 43 - Oakland
$ws.Range("W43").Value = 1

# Row 47 - Little Italy, UIC
$ws.Range("AD47").Value = 3

# Row 49 - Grand Boulevard
$ws.Range("B49").Value = 2

# Row 50 - Garfield Ridge
$ws.Range("W50").Value = 1

# Row 52 - Chatham
$ws.Range("B52").Value = 3

# Row 53 - Calumet Heights
$ws.Range("I53").Value = 3

# Row 59 - Archer Heights
$ws.Range("B59").Value = 1

# Row 60 - Armour Square
$ws.Range("P60").Value = 1

# Row 65 - Chicago Lawn
$ws.Range("AK65").Value = 2

# Row 78 - Lake View
$ws.Range("AD78").Value = 1
